$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '27.126.53'
    'E2' = '  +3.31%  '
    'D3' = '1.656.58'
    'E3' = '  +3.56%  '
    'E4' = '  -0.21%  '
    'D5' = '215.31'
    'E5' = '  +1.24%  '
    'E6' = '  +1.38%  '
    'E7' = '  -0.25%  '
    'E8' = '  +2.12%  '
    'D9' = '0.0615'
    'E9' = '  +1.40%  '
    'D10' = '19.52'
    'E10' = '  +2.94%  '
    'D11' = '0.0862'
    'E11' = '  +0.60%  '
    'D12' = '1.889.28'
    'E12' = '  +3.53%  '
    'D13' = '1.674.60'
    'E13' = '  +4.49%  '
    'D14' = '4.08'
    'E14' = '  +1.71%  '
    'E15' = '  +2.87%  '
    'E16' = '  +1.94%  '
    'D17' = '240.97'
    'E17' = '  +5.31%  '
    'D18' = '27.094.66'
    'E18' = '  +3.14%  '
    'E19' = '  +3.09%  '
    'D20' = '0.0₃0729'
    'E20' = '  +1.41%  '
    'D21' = '0.999'
    'E21' = '  -0.19%  '
    'E22' = '  +4.34%  '
    'D23' = '2.25'
    'E23' = '  +4.10%  '
    'D24' = '9.28'
    'E24' = '  +3.65%  '
    'D25' = '146.18'
    'E25' = '  +0.41%  '
    'E26' = '  -0.21%  '
    'E27' = '  +2.42%  '
    'E28' = '  +1.28%  '
    'D29' = '15.87'
    'E29' = '  +2.93%  '
    'D30' = '0.0497'
    'E30' = '  +0.67%  '
    'E31' = '  +0.71%  '
    'D32' = '1.522.89'
    'E32' = '  +5.40%  '
    'E33' = '  +2.70%  '
    'D34' = '3.05'
    'E34' = '  +3.15%  '
    'D35' = '1.57'
    'E35' = '  +6.94%  '
    'D36' = '2.43'
    'E36' = '  -0.17%  '
    'E37' = '  +1.43%  '
    'D38' = '0.901'
    'E38' = '  +9.79%  '
    'E39' = '  +2.61%  '
    'E40' = '  +3.12%  '
    'E41' = '  -0.28%  '
    'D42' = '2.28'
    'E42' = '  +4.77%  '
    'D43' = '65.14'
    'E43' = '  +7.39%  '
    'D44' = '1.795.77'
    'E44' = '  +3.32%  '
    'D45' = '0.773'
    'E45' = '  +1.66%  '
    'D46' = '0.917'
    'E46' = '  -0.81%  '
    'D47' = '90.54'
    'E47' = '  +3.50%  '
    'E48' = '  +3.32%  '
    'D49' = '0.0₆0102'
    'E50' = '  +0.54%  '
    'D51' = '0.0977'
    'E51' = '  +2.78%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
